$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Item("2022-Q1")

# ---------------------------------------------------------------------------
# 1. Duplicate the existing "2022-Q1" sheet, placing the copy right after it.
#    The copy retains all the original Q1 fund-holdings data & formatting
#    and becomes the new, separate "2022-Q1" sheet.
# ---------------------------------------------------------------------------
$q1.Copy($null, $q1)
$q1Copy = $wb.Worksheets.Item($q1.Index + 1)
$q1Copy.Name = "2022-Q1_tmp"

# 2. The original sheet (still holding Q1 data for now) becomes the Q3 sheet.
$q1.Name = "2022-Q3"
$q3 = $q1

# 3. Rename the duplicate back to "2022-Q1".
$q1Copy.Name = "2022-Q1"

# ---------------------------------------------------------------------------
# 4. Replace Q3's old content with the new quarter's fund-holdings data.
# ---------------------------------------------------------------------------
$q3.Cells.Clear()

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q3.Cells.Item(1, 2 + $i).Value = $headers[$i]
}

$rows = @(
    @("003359","大成中证360互联网+大数据100指数C","1.11","92.17","1.07","0.0119",3),
    @("002236","大成中证360互联网+大数据100指数A","1.03","92.17","1.07","0.0110",3),
    @("005536","渤海汇金量化成长混合","0.38","84.80","1.58","0.0060",8),
    @("010584","渤海汇金新动能主题混合","0.27","91.08","1.92","0.0052",6),
    @("000804","中信建投稳利混合A","0.21","38.27","1.75","0.0037",9),
    @("006844","中信建投稳利混合C","0.10","38.27","1.75","0.0018",9)
)

# Force column B:G to text BEFORE writing, so numeric-looking values (fund
# codes, percentages, …) stay text instead of being coerced to numbers.
$q3.Range("B2:G7").NumberFormat = "@"

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    $rowIdx = 2 + $r
    $q3.Cells.Item($rowIdx, 1).Value = $r
    for ($c = 0; $c -lt 6; $c++) {
        $q3.Cells.Item($rowIdx, 2 + $c).Value = $row[$c]
    }
    $q3.Cells.Item($rowIdx, 8).Value = $row[6]
}

# Match the workbook's existing header/index-column styling (bold, bordered,
# centered) by copying the format already used on the "总计" sheet.
$total.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$total.Range("A2").Copy()
$q3.Range("A2:A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Output "done"
